$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "- КАТЕГОРИЯ ВОДИТЕЛЬСКОГО УДОСТОВЕРЕНИЯ (...)"
# The original paragraph holds the bullet text split across three runs:
#   <w:tab/><w:t>– КАТЕГОРИ</w:t>  +  <w:t>Я</w:t>  +  <w:t> ВОДИТЕЛЬСКОГО ...</w:t>
# After the edit they collapse into a single run (the leading <w:tab/> stays
# a distinct element, it must not be swallowed into the text run).
# ---------------------------------------------------------------------------
$categoryParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "*КАТЕГОРИЯ ВОДИТЕЛЬСКОГО УДОСТОВЕРЕНИЯ*") {
        $categoryParaIndex = $i
    }
}

$p1 = $d.Paragraphs.Item($categoryParaIndex)
$r1 = $p1.Range
$paraStart1 = $r1.Start
$paraEnd1 = $r1.End
$runsRange1 = $d.Range($paraStart1, $paraEnd1 - 1)

$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="28"/><w:lang w:val="zxx" w:eastAsia="zxx" w:bidi="zxx"/></w:rPr><w:tab/><w:t>– КАТЕГОРИЯ ВОДИТЕЛЬСКОГО УДОСТОВЕРЕНИЯ (ВОДИТЕЛЬСКОЕ УДОСТОВЕРЕНИЕ, КАТЕГОРИЯ) (M:N);</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$runsRange1.InsertXML($xmlFrag1)

# ---------------------------------------------------------------------------
# Change 2: the closing "P.S." paragraph.
# Its second run (everything after the italic "P.S." run) gets its text
# extended and split into five runs that share the same formatting.
# ---------------------------------------------------------------------------
$psParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "*в дальнейшей работе допускается сокращение понятия*") {
        $psParaIndex = $i
    }
}

$p2 = $d.Paragraphs.Item($psParaIndex)
$r2 = $p2.Range
$text2 = $r2.Text
$relIdx2 = $text2.IndexOf(": в дальнейшей")
$run2Start = $r2.Start + $relIdx2
$paraEnd2 = $r2.End
$runsRange2 = $d.Range($run2Start, $paraEnd2 - 1)

$rPrCommon = '<w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="28"/><w:u w:val="none"/><w:lang w:val="zxx" w:eastAsia="zxx" w:bidi="zxx"/></w:rPr>'

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
  '<w:r>' + $rPrCommon + '<w:t>: в дальнейшей работе допускается сокращение понятия "ТРАНСПОРТНОЕ СРЕДСТВО" до аббревиатуры "ТС", поняти</w:t></w:r>' + `
  '<w:r>' + $rPrCommon + '<w:t>я</w:t></w:r>' + `
  '<w:r>' + $rPrCommon + '<w:t xml:space="preserve"> "ВОДИТЕЛЬСКОЕ УДОСТОВЕРЕНИЕ"</w:t></w:r>' + `
  '<w:r>' + $rPrCommon + '<w:t xml:space="preserve"> –</w:t></w:r>' + `
  '<w:r>' + $rPrCommon + '<w:t xml:space="preserve"> до аббревиатуры "ВУ".</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$runsRange2.InsertXML($xmlFrag2)
